$wb = $excel.ActiveWorkbook

# --- Add the new "cronograma" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "cronograma"

# --- Header row (row 1) ---
$ws.Range("A1").Value = "idestimacion"
$ws.Range("B1").Value = " idtarea"
$ws.Range("C1").Value = " incluir"
$ws.Range("D1").Value = " porcentaje"
$ws.Range("E1").Value = " recursos"
$ws.Range("F1").Value = " dias"
$ws.Range("G1").Value = " horas"

# --- Data rows ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 40

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.3
$ws.Range("E3").Value = 1

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0.2
$ws.Range("E4").Value = 1

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0.4
$ws.Range("E5").Value = 1

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0.3
$ws.Range("E6").Value = 1

$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0.1
$ws.Range("E7").Value = 1

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 1
$ws.Range("E8").Value = 1

# --- Underlined-font style used for G3 and the H:I "extra" columns ---
$ws.Range("G3").Font.Underline = $true

# Touch I8 momentarily so the used range extends through column I / row 8,
# then turn the underline back off so the cell itself stays blank/unstyled.
$extra = $ws.Range("I8")
$extra.Font.Underline = $true
$extra.Font.Underline = $false

# --- Column widths ---
$ws.Range("A1:G1").EntireColumn.ColumnWidth = 16.8
$ws.Range("H1:I1").EntireColumn.ColumnWidth = 9.59

# --- Selection / active cell ---
$ws.Range("D4").Select() | Out-Null
